$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44939, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Candy White", "Primera", 10, 370000, 380000, 375000, "`$/bins (420 kilos)", "Región de O'Higgins", 893, 420),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44939, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Candy White", "Segunda", 10, 320000, 330000, 325000, "`$/bins (420 kilos)", "Región de O'Higgins", 774, 420),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44939, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Super Queen", "Especial", 10, 430000, 440000, 435000, "`$/bins (420 kilos)", "Región de O'Higgins", 1036, 420),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44939, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Super Queen", "Primera", 16, 380000, 390000, 385000, "`$/bins (420 kilos)", "Región de O'Higgins", 917, 420),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44939, 4, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Super Queen", "Segunda", 10, 340000, 350000, 345000, "`$/bins (420 kilos)", "Región de O'Higgins", 821, 420)
)

$startRow = 414
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
